$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 405.7
$ws.Range("I33").Value = 442.75
$ws.Range("J33").Value = 257.5
$ws.Range("K33").Value = 442.75
$ws.Range("L33").Value = 257.5
$ws.Range("M33").Value = -213.75
$ws.Range("N33").Value = -715.5

$ws.Range("H43").Value = 291666900
$ws.Range("I43").Value = 1000000000
$ws.Range("J43").Value = 55555892
$ws.Range("K43").Value = 1000000000
$ws.Range("L43").Value = 55555892
$ws.Range("M43").Value = -999999931
$ws.Range("N43").Value = -55556030

$ws.Range("H57").Value = 12784
$ws.Range("J57").Value = 12784
$ws.Range("L57").Value = 38352
$ws.Range("N57").Value = -39350

$ws.Range("H74").Value = 3838.75
$ws.Range("I74").Value = 3530
$ws.Range("K74").Value = 3530
$ws.Range("M74").Value = -2594

$ws.Range("H77").Value = 3838.75
$ws.Range("I77").Value = 3530
$ws.Range("K77").Value = 17650
$ws.Range("M77").Value = -12970

$ws.Range("H113").Value = 3404.5454
$ws.Range("I113").Value = 2862.5
$ws.Range("J113").Value = 4850
$ws.Range("K113").Value = 2862.5
$ws.Range("L113").Value = 4850
$ws.Range("M113").Value = 391.5
$ws.Range("N113").Value = -11358

$ws.Range("H116").Value = 2257.353
$ws.Range("J116").Value = 3525.75
$ws.Range("L116").Value = 3525.75
$ws.Range("N116").Value = -10409.75

$ws.Range("H129").Value = 22125.213
$ws.Range("I129").Value = 338.7143
$ws.Range("J129").Value = 25937.85
$ws.Range("K129").Value = 1016.1429
$ws.Range("L129").Value = 77813.54999999999
$ws.Range("M129").Value = 3983.8571
$ws.Range("N129").Value = -87813.54999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20542

$ws.Range("H43").Value = 7125.6665
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 7125.6665
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 7125.6665
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -7751.6665

$ws.Range("H63").Value = 2001760
$ws.Range("I63").Value = 5000900
$ws.Range("K63").Value = 5000900
$ws.Range("M63").Value = -5000214

$ws.Range("H66").Value = 2001760
$ws.Range("I66").Value = 5000900
$ws.Range("K66").Value = 25004500
$ws.Range("M66").Value = -25001068

$ws.Range("H102").Value = 1823.5
$ws.Range("I102").Value = 1980
$ws.Range("K102").Value = 1980
$ws.Range("M102").Value = -358

$ws.Range("H122").Value = 1618
$ws.Range("I122").Value = 1481.75
$ws.Range("J122").Value = 1799.6666
$ws.Range("K122").Value = 4445.25
$ws.Range("L122").Value = 5398.9998
$ws.Range("M122").Value = -1995.25
$ws.Range("N122").Value = -10298.9998

$ws.Range("H132").Value = 7057.84
$ws.Range("I132").Value = 8934.25
$ws.Range("J132").Value = 3722
$ws.Range("K132").Value = 26802.75
$ws.Range("L132").Value = 11166
$ws.Range("M132").Value = -24272.75
$ws.Range("N132").Value = -16226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1486.5238
$ws.Range("I99").Value = 1433.1818
$ws.Range("J99").Value = 1545.2
$ws.Range("K99").Value = 1433.1818
$ws.Range("L99").Value = 1545.2
$ws.Range("M99").Value = 64.81819999999993
$ws.Range("N99").Value = -4541.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 3574134
$ws.Range("I132").Value = 2393.2942
$ws.Range("J132").Value = 6947445
$ws.Range("K132").Value = 7179.882599999999
$ws.Range("L132").Value = 20842335
$ws.Range("M132").Value = -4649.882599999999
$ws.Range("N132").Value = -20847395

$ws.Range("H134").Value = 1153.5
$ws.Range("I134").Value = 1138
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 3414
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -879
$ws.Range("N134").Value = -8670

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3011.111
$ws.Range("I81").Value = 1800
$ws.Range("J81").Value = 3162.5
$ws.Range("K81").Value = 5400
$ws.Range("L81").Value = 9487.5
$ws.Range("M81").Value = -4277
$ws.Range("N81").Value = -11733.5

$ws.Range("H84").Value = 3011.111
$ws.Range("I84").Value = 1800
$ws.Range("J84").Value = 3162.5
$ws.Range("K84").Value = 16200
$ws.Range("L84").Value = 28462.5
$ws.Range("M84").Value = -10584
$ws.Range("N84").Value = -39694.5

$ws.Range("H131").Value = 758.58
$ws.Range("J131").Value = 791.3913
$ws.Range("L131").Value = 2374.1739
$ws.Range("N131").Value = -12454.1739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2256.2856
$ws.Range("I126").Value = 3098.2856
$ws.Range("J126").Value = 1414.2858
$ws.Range("K126").Value = 9294.856800000001
$ws.Range("L126").Value = 4242.857400000001
$ws.Range("M126").Value = -6824.856800000001
$ws.Range("N126").Value = -9182.857400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 695.2
$ws.Range("I22").Value = 695.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 695.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -400.2
$ws.Range("N22").Value = ""

$ws.Range("H27").Value = 695.2
$ws.Range("I27").Value = 695.2
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 695.2
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -588.2
$ws.Range("N27").Value = ""

$ws.Range("H100").Value = 33334880
$ws.Range("I100").Value = 66667900
$ws.Range("J100").Value = 1860
$ws.Range("K100").Value = 66667900
$ws.Range("L100").Value = 1860
$ws.Range("M100").Value = -66667359
$ws.Range("N100").Value = -2942

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 419.09525
$ws.Range("I113").Value = 454.44446
$ws.Range("J113").Value = 392.58334
$ws.Range("K113").Value = 1363.33338
$ws.Range("L113").Value = 1177.75002
$ws.Range("M113").Value = 806.66662
$ws.Range("N113").Value = -5517.750019999999
